$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper cells we reuse as format donors (already present in the workbook):
#   E1   -> "DONE" style (fillId 2 / centered)               -> s="1"
#   D3   -> wrapped / colored objective-row style             -> s="6"
#   C14  -> centered style                                    -> s="2"
# ---------------------------------------------------------------------------

# --- D6:D12 pick up the "objective" style already used by D3:D5 ------------
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null

# --- new row 12 "Implement Jump to bootloader command in python" -----------
$ws.Range("D12").Value = "Implement Jump to bootloader command in python"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

# --- new "DONE" markers in column E for rows 6-14 ---------------------------
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Value = "DONE"
$ws.Range("E7").Value = "DONE"
$ws.Range("E8").Value = "DONE"
$ws.Range("E9").Value = "DONE"
$ws.Range("E10").Value = "DONE"
$ws.Range("E11").Value = "DONE"
$ws.Range("E12").Value = "DONE"
$ws.Range("E14").Value = "DONE"

# --- new "DONE" markers in column H for rows 11,12,14 -----------------------
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null

$ws.Range("H11").Value = "DONE"
$ws.Range("H12").Value = "DONE"
$ws.Range("H14").Value = "DONE"

# --- bootloader command table re-shuffle (rows 11-14, cols I/J/K) ----------
# Row 12 now shows the 0xB4 / Flash Application pairing, row 13 shows the
# 0xB5 / Check Application Integrity pairing (+ its description), and row 14
# gains the new 0xB6 / Jump to Application pairing.
$ws.Range("J12").Value = "Flash Application"
$ws.Range("K12").ClearContents()

$ws.Range("J13").Value = "Check Application Integrity"
$ws.Range("K13").Value = "Calculate the CRC of Application and Verify it with CRC already STORED"

$ws.Range("I14").Value = "0xB6"
$ws.Range("J14").Value = "Jump to Application"

# ---------------------------------------------------------------------------
# New packet-format table (rows 18-30)
# ---------------------------------------------------------------------------
$ws.Range("G18").Value = 1
$ws.Range("I18").Value = 4
$ws.Range("C14").Copy() | Out-Null
$ws.Range("G18").PasteSpecial(-4122) | Out-Null
$ws.Range("G18").Value = 1
$ws.Range("C14").Copy() | Out-Null
$ws.Range("I18").PasteSpecial(-4122) | Out-Null
$ws.Range("I18").Value = 4

$cmdRows = 19, 21, 23, 25, 27, 29
$headers = "NO.OF.PACKTES - 2Bytes", "LEN.NEXT.PACK - 2BYTES", "PACKET", "LEN.NEXT.PACK", "PACKET", "END.PACK"
for ($i = 0; $i -lt $cmdRows.Length; $i++) {
    $r = $cmdRows[$i]
    $ackRow = $r + 1
    $ws.Range("G$r").Value = "CMD"
    $ws.Range("H$r").Value = $headers[$i]
    $ws.Range("I$r").Value = "CRC"

    if ($r -ne 29) {
        $ws.Range("H9").Copy() | Out-Null
        $ws.Range("J$r").PasteSpecial(-4122) | Out-Null
        $ws.Range("J$r").Value = "DONE"
    }

    if ($ackRow -le 28) {
        $ws.Range("G$ackRow").Value = "ACK"
        $ws.Range("H9").Copy() | Out-Null
        $ws.Range("J$ackRow").PasteSpecial(-4122) | Out-Null
        $ws.Range("J$ackRow").Value = "DONE"
    }
}

# --- column H sizing (bestFit-like autosize from the short header labels) --
# Done before H30 gets its long wrapped paragraph so that cell doesn't blow
# out the computed best-fit width.
$ws.Columns("H").AutoFit() | Out-Null

# --- closing "End pack" description, wrapped, in H30 ------------------------
# (WrapText is applied first so the new "wrap only" style is minted before
# the red-fill style used by J29, matching the original authoring order.)
$ws.Range("H30").Value = "End pack is length of 4 bytes with CRC of complete BIN file, `nSimilaryly STM also calculates the CRC of whole data in PFLASH and sends it to Python.`nPython Compares the CRC with calculated CRC and confirms if Application is flashed successfully."
$ws.Range("H30").WrapText = $true
$ws.Rows(30).RowHeight = 187.2

# Row 29's ACK lives on row 29 itself per the diff (G29/H29/I29), and J29 is
# an empty, red-filled cell rather than another DONE marker.
$ws.Range("J29").Interior.Color = 255

# --- selection matches the post-edit workbook state -------------------------
$ws.Range("H5").Select() | Out-Null
